$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(34352, "Vitor Hugo das Neves", "Vendas", "Viagem de negocios", 5, 45094, 7204),
    @(16397, "Sophia Moraes", "Marketing", "Outros", 2, 45100, 5523.55),
    @(36677, "Laís Fogaça", "Vendas", "Viagem de negocios", 6, 45080, 2535.95),
    @(9764, "Cecilia Vieira", "Recursos Humanos", "Doenca", 1, 45101, 7617.81),
    @(59442, "Ana Carolina da Mata", "P&D", "Consulta medica", 1, 45102, 6352.33),
    @(9068, "Rebeca Pastor", "TI", "Doenca", 5, 45088, 9030.99),
    @(81250, "Diogo da Costa", "Financeiro", "Outros", 4, 45096, 6671.58),
    @(410, "Dr. Francisco Montenegro", "Engenharia", "Outros", 5, 45088, 9745.190000000001),
    @(35113, "Caleb Cunha", "Operacoes", "Consulta medica", 1, 45095, 9158.9),
    @(77887, "Lucas Gabriel Gonçalves", "Juridico", "Outros", 3, 45104, 2971.69)
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    $ws.Cells.Item($row, 5).Value = $record[4]
    $ws.Cells.Item($row, 6).Value = $record[5]
    $ws.Cells.Item($row, 7).Value = $record[6]
    $row++
}
